# Update the two "Valor Mora" figures that were swapped between the
# period-2002 row (16) and the period-1908 row (22) on Hoja1.
#   F16: 18771 -> 33125
#   F22: 33125 -> 18771
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$ws.Range("F16").Value = 33125
$ws.Range("F22").Value = 18771
